# Weekly refresh of the "Cebollín baby" price sheet: a new weekly record is
# inserted as row 32 (pushing the previously-existing rows 32-44 down to
# 33-45, unchanged), and the new row is populated with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 32 - everything below (old rows 32-44)
# shifts down to 33-45, carrying its values and formatting (e.g. the date
# style on column D) along with it.
$ws.Rows(32).Insert()

# Populate the newly inserted row 32 with the new weekly observation.
$ws.Cells.Item(32, 1).Value  = 1
$ws.Cells.Item(32, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(32, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(32, 4).Value  = 44452
$ws.Cells.Item(32, 5).Value  = 15
$ws.Cells.Item(32, 6).Value  = 100112038
$ws.Cells.Item(32, 7).Value  = "Cebollín baby"
$ws.Cells.Item(32, 8).Value  = "Sin especificar"
$ws.Cells.Item(32, 9).Value  = "Primera"
$ws.Cells.Item(32, 10).Value = 300
$ws.Cells.Item(32, 11).Value = 1900
$ws.Cells.Item(32, 12).Value = 2000
$ws.Cells.Item(32, 13).Value = 1950
$ws.Cells.Item(32, 14).Value = "$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(32, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(32, 16).Value = 975
$ws.Cells.Item(32, 17).Value = 2
$ws.Cells.Item(32, 18).Value = "Hortaliza"
